$wb = $excel.ActiveWorkbook

$wsBOM = $wb.Worksheets.Item("BOM Report")
$wsInfo = $wb.Worksheets.Item("Project Information")

# --- Update generated/report date-time (regenerated BOM report) ---
# Leading apostrophe forces these to stay plain text (matching the
# original quotePrefix-styled cells) instead of being auto-converted
# to date/time serial values.
# "Generated On:" (BOM Report!B6) and "Report Date & Time" (Project
# Information!B10) both held the same "<date> <time>" string, so both
# are updated together to the new report timestamp.
$wsBOM.Range("B6").Value = "'2020-01-20 9:50 PM"
$wsInfo.Range("B10").Value = "'2020-01-20 9:50 PM"
$wsInfo.Range("B8").Value = "'9:50 PM"
$wsInfo.Range("B9").Value = "'2020-01-20"

# --- Update unit price (column G) and subtotal (column I) values ---
# Row 12
$wsBOM.Range("G12").Value = 0.15661
# Row 13
$wsBOM.Range("G13").Value = 0.19575999999999999
# Row 14
$wsBOM.Range("G14").Value = 0.29232999999999998
# Row 15
$wsBOM.Range("G15").Value = 0.26101000000000002
# Row 16
$wsBOM.Range("G16").Value = 0.18271000000000001
$wsBOM.Range("I16").Value = 0.18271000000000001
# Row 17
$wsBOM.Range("G17").Value = 0.18271000000000001
# Row 23
$wsBOM.Range("I23").Value = 5.25
# Row 25
$wsBOM.Range("G25").Value = 0.030020000000000002
$wsBOM.Range("I25").Value = 0.45023999999999997
# Row 26
$wsBOM.Range("G26").Value = 0.13050999999999999
$wsBOM.Range("I26").Value = 0.91354000000000002
# Row 27
$wsBOM.Range("G27").Value = 0.074389999999999998
# Row 28
$wsBOM.Range("G28").Value = 0.030020000000000002
$wsBOM.Range("I28").Value = 0.33017999999999997
# Row 29
$wsBOM.Range("G29").Value = 0.13050999999999999
$wsBOM.Range("I29").Value = 0.91354000000000002
# Row 30
$wsBOM.Range("G30").Value = 0.31320999999999999
$wsBOM.Range("I30").Value = 0.31320999999999999
# Row 32
$wsBOM.Range("I32").Value = 6.79
# Row 33
$wsBOM.Range("G33").Value = 1.74

$wb.Application.Calculate()
